$wb = $excel.ActiveWorkbook

# "Overview" sheet: Latest HO Xliff Generate Date for 6edb2551-... row
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-17 22:44:34"

# "zh-cn" sheet: Correspond Handoff Datetime / Correspond Handback DateTime for 6edb2551-... row
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-08-17 22:44:28"
$wsZhCn.Range("K4").Value = "2016-08-17 22:44:57"

# "de-de" sheet: Correspond Handoff Datetime / Correspond Handback DateTime for 6edb2551-... row
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H4").Value = "2016-08-17 22:44:34"
$wsDeDe.Range("K4").Value = "2016-08-17 22:45:12"
